# Fill in missing/corrected value for the primary dataset using cleaned
# value from a different dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 corresponds to "1937" (Chinese Population %) - corrected from 56.1 to 55.1
$ws.Range("B8").Value = 55.1

$wb.Save()
